$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the "custom accuracy" (2-decimal) rounded figures
$ws.Range("B5").Value = 10.57
$ws.Range("C5").Value = 7.61
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 22.96
$ws.Range("F5").Value = 18.43
$ws.Range("G5").Value = 8.289999999999999
$ws.Range("H5").Value = 32.67
$ws.Range("I5").Value = 12.8
$ws.Range("J5").Value = 5.58
$ws.Range("K5").Value = 8.17
$ws.Range("L5").Value = 9.210000000000001
$ws.Range("M5").Value = 9.67
$ws.Range("N5").Value = 2.66
$ws.Range("O5").Value = 8.27
$ws.Range("P5").Value = 11.68
$ws.Range("Q5").Value = 7.13
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.55
$ws.Range("T5").Value = 118.57
$ws.Range("U5").Value = 23.14
$ws.Range("V5").Value = 7.63
$ws.Range("W5").Value = 15.37
$ws.Range("X5").Value = 8.039999999999999
$ws.Range("Y5").Value = 1.42
$ws.Range("Z5").Value = 15.72
$ws.Range("AA5").Value = 6.74
$ws.Range("AB5").Value = 6.05
$ws.Range("AC5").Value = 7.12
$ws.Range("AD5").Value = 9.66
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 29.76
$ws.Range("AG5").Value = 4.21
$ws.Range("AH5").Value = 9.550000000000001

# Remove row 6 entirely (data trimmed), which also updates the sheet dimension
$ws.Rows.Item(6).Delete()

$wb.Save()
